$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: target cell -> new literal text for that cell.
# 'ForceText' entries are values that look like plain decimal numbers
# (e.g. "1.005"); those are written with a leading apostrophe so Excel
# keeps storing them as literal text (matching the inlineStr source data)
# instead of silently coercing them to a Double, and the quote-prefix
# cell style that the apostrophe trick implies is reset right back to
# "Normal" so no stray formatting is introduced.
$updates = @(
    [pscustomobject]@{ Cell = 'D2'; Value = '23.805.13'; ForceText = $false }
    [pscustomobject]@{ Cell = 'E2'; Value = '  -2.63%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D3'; Value = '1.617.08'; ForceText = $false }
    [pscustomobject]@{ Cell = 'E3'; Value = '  -3.17%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D4'; Value = '1.005'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E4'; Value = '  -0.08%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D6'; Value = '306.44'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E6'; Value = '  -2.12%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D7'; Value = '0.3899'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E7'; Value = '  +0.07%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D8'; Value = '0.3832'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E8'; Value = '  -2.31%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D9'; Value = '1.006'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E9'; Value = '  -0.07%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D10'; Value = '1.348'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E10'; Value = '  -3.58%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D11'; Value = '48.77'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E11'; Value = '  -6.04%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D12'; Value = '0.08393'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E12'; Value = '  -2.33%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D13'; Value = '23.65'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E13'; Value = '  -4.97%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D14'; Value = '6.951'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E14'; Value = '  -4.35%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D15'; Value = '0.00001268'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E15'; Value = '  -2.99%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D16'; Value = '7.397'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E16'; Value = '  -3.90%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D17'; Value = '1.616.77'; ForceText = $false }
    [pscustomobject]@{ Cell = 'E17'; Value = '  -3.56%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D18'; Value = '93.01'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E18'; Value = '  -0.34%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D19'; Value = '0.06917'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E19'; Value = '  -1.92%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D20'; Value = '19.74'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E20'; Value = '  -2.89%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D21'; Value = '6.778'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E21'; Value = '  -3.56%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D22'; Value = '1.003'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E22'; Value = '  -0.29%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D23'; Value = '13.34'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E23'; Value = '  -4.00%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D24'; Value = '23.811.79'; ForceText = $false }
    [pscustomobject]@{ Cell = 'E24'; Value = '  -2.62%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D25'; Value = '2.407'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E25'; Value = '  +1.37%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D26'; Value = '2.761'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E26'; Value = '  +1.22%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D27'; Value = '22.04'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E27'; Value = '  -4.97%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D28'; Value = '157.22'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E28'; Value = '  -2.40%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D29'; Value = '139.24'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E29'; Value = '  -5.21%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D30'; Value = '5.268'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E30'; Value = '  -9.04%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D31'; Value = '7.784'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E31'; Value = '  -5.01%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D32'; Value = '2.466'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E32'; Value = '  -2.04%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D33'; Value = '1.791.76'; ForceText = $false }
    [pscustomobject]@{ Cell = 'E33'; Value = '  -3.89%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D34'; Value = '0.08006'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E34'; Value = '  -3.55%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'B35'; Value = 'ImmutableX'; ForceText = $false }
    [pscustomobject]@{ Cell = 'C35'; Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; ForceText = $false }
    [pscustomobject]@{ Cell = 'D35'; Value = '0.9644'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E35'; Value = '  +0.00%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'B36'; Value = 'VeChain'; ForceText = $false }
    [pscustomobject]@{ Cell = 'C36'; Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; ForceText = $false }
    [pscustomobject]@{ Cell = 'D36'; Value = '0.02863'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E36'; Value = '  -4.75%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D37'; Value = '6.589'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E37'; Value = '  -4.94%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D38'; Value = '0.2645'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E38'; Value = '  -5.04%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D39'; Value = '0.09103'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E39'; Value = '  -3.76%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D40'; Value = '10.34'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E40'; Value = '  +1.43%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D41'; Value = '13.29'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E41'; Value = '  -1.06%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D42'; Value = '1.413'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E42'; Value = '  -6.44%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D43'; Value = '0.7425'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E43'; Value = '  -4.97%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D44'; Value = '15.78'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E44'; Value = '  -3.36%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D45'; Value = '0.6822'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E45'; Value = '  -3.16%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D46'; Value = '2.427'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E46'; Value = '  -4.17%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D47'; Value = '4.051'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E47'; Value = '  -2.87%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'E48'; Value = '  +0.04%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D49'; Value = '0.08218'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E49'; Value = '  -4.16%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D50'; Value = '132.27'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E50'; Value = '  -3.35%  '; ForceText = $false }
    [pscustomobject]@{ Cell = 'D51'; Value = '1.243'; ForceText = $true }
    [pscustomobject]@{ Cell = 'E51'; Value = '  -5.25%  '; ForceText = $false }
)

foreach ($update in $updates) {
    $range = $ws.Range($update.Cell)
    if ($update.ForceText) {
        # Leading apostrophe forces text storage for number-looking strings.
        $range.Value = "'" + $update.Value
        $range.Style = "Normal"
    } else {
        $range.Value = $update.Value
    }
}
